# Regenerate merged AHB files
# - Rename the "*_old" / "*_new" header columns to "*_FV2404" / "*_FV2410"
# - Wrap the data range in an Excel Table ("Table1")
# - Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row labels: *_old -> *_FV2404, *_new -> *_FV2410
$ws.Cells.Replace("_old", "_FV2404")
$ws.Cells.Replace("_new", "_FV2410")

# 2) Turn the A1:U79 range into a native Excel table named Table1
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U79"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# 3) Freeze the top (header) row
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
